$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated training schedule values (rows 2-5) and new row 6
$data = @(
    @(1, 0, 3, 4, 5, 4, 2, 54, 5),
    @(2, 1, 2, 6, 3, 5, 1, 65, 5),
    @(3, 1, 4, 2, 9, 1, 5, 21, 5),
    @(4, 0, 1, 3, 4, 3, 3, 43, 5),
    @(5, 3, 4, 5, 8, 2, 4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowVals[$col - 1]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$null = $ws.Range("I1").Select()
